$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (shared string changes) ---
$ws.Range("B1").Value = "Severe Storm(s)"
$ws.Range("I1").Value = "Volcanic Eruption"
$ws.Range("J1").Value = "Severe Storm"
$ws.Range("O1").Value = "Snowstorm"

# --- Clear column B for rows 2-70 (years 1953-2021) ---
$ws.Range("B2:B70").ClearContents()

# --- Update row 71 (year 2022) values ---
$ws.Range("B71").Value = 3
$ws.Range("C71").Value = 41
$ws.Range("D71").Value = 10
$ws.Range("G71").Value = 10
$ws.Range("H71").Value = 2
$ws.Range("J71").Value = 24
$ws.Range("R71").Value = 2

# --- Add new row 72 (year 2023) ---
$ws.Range("A72").Value = 2023
$ws.Range("C72").Value = 1
$ws.Range("D72").Value = 2
$ws.Range("G72").Value = 3
$ws.Range("J72").Value = 1
